$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Installation")

# Row 79: same style as B78 (applied via copy/paste of formats), new shared string
$ws1.Range("B78").Copy()
$ws1.Range("B79").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws1.Range("B79").Value = "npm i vee-validate --save"

# Row 81: plain new cell (no style)
$ws1.Range("B81").Value = " npm install -S yup   "

# Row 83: plain new cell (no style)
$ws1.Range("B83").Value = "npm install pinia"

# Update the view: active cell and scroll position
$ws1.Activate()
$ws1.Range("H82").Select()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
